$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.780.78"
$ws.Range("E2").Value = "  -2.63%  "
$ws.Range("D3").Value = "3.487.00"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.70"
$ws.Range("E5").Value = "  -2.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.36"
$ws.Range("E6").Value = "  -4.17%  "
$ws.Range("D7").Value = "3.488.11"
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.380"
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").Value = "4.077.90"
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.37"
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000177"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").Value = "3.458.18"
$ws.Range("E17").Value = "  -3.76%  "
$ws.Range("D18").Value = "63.895.47"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.06"
$ws.Range("E19").Value = "  +2.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.36"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "384.04"
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.575"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").Value = "3.628.43"
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.03"
$ws.Range("E25").Value = "  -2.58%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000111"
$ws.Range("E27").Value = "  -3.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.57"
$ws.Range("E28").Value = "  -3.94%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.40"
$ws.Range("E30").Value = "  -6.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.23"
$ws.Range("E31").Value = "  -3.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.18"
$ws.Range("E32").Value = "  -4.08%  "
$ws.Range("D33").Value = "3.495.94"
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.62"
$ws.Range("E35").Value = "  -2.70%  "
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.21"
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.57"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0796"
$ws.Range("E41").Value = "  -3.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.82"
$ws.Range("E42").Value = "  +2.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.807"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.47"
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("E46").Value = "  -4.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.37"
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.63"
$ws.Range("E48").Value = "  -3.03%  "
$ws.Range("D49").Value = "2.434.01"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.83"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.885"
$ws.Range("E51").Value = "  -1.44%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.86"
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.55"
$ws.Range("E39").Value = "  -2.70%  "
